$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "NC"

$ws2.Range("B1").Value = "In-vehicle"
$ws2.Range("C1").Value = "At-stop"
$ws2.Range("D1").Value = "Extra"
$ws2.Range("E1").Value = "Total"
$ws2.Range("A2").Value = "No control"
$ws2.Range("B2").Value = 2101.086661275402
$ws2.Range("C2").Value = 12498.70440518066
$ws2.Range("D2").Value = 141.4698672425732
$ws2.Range("E2").Value = 14741.26093369863

$ws1.Range("B1:E1").Copy()
$ws2.Range("B1:E1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

[void]($excel.CutCopyMode = $false)

[void]$ws2.Range("A1").Select()
[void]$ws1.Select()
